$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.852.36"
$ws.Range("E2").Value = "  +0.97%  "
# Row 3
$ws.Range("D3").Value = "1.768.13"
$ws.Range("E3").Value = "  +0.86%  "
# Row 4
$ws.Range("D4").Value = "'1.006"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.52%  "
# Row 5
$ws.Range("D5").Value = "'327.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.95%  "
# Row 6
$ws.Range("D6").Value = "'1.004"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.37%  "
# Row 7
$ws.Range("D7").Value = "'0.4485"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.03%  "
# Row 8
$ws.Range("D8").Value = "'0.3575"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.61%  "
# Row 9
$ws.Range("B9").Value = "OKB"
$ws.Range("C9").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D9").Value = "'42.31"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.39%  "
# Row 10
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Value = "'0.07457"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.39%  "
# Row 11
$ws.Range("D11").Value = "'1.096"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.22%  "
# Row 12
$ws.Range("D12").Value = "'1.006"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.49%  "
# Row 13
$ws.Range("D13").Value = "'20.88"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.55%  "
# Row 14
$ws.Range("D14").Value = "'6.050"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.81%  "
# Row 15
$ws.Range("D15").Value = "'7.210"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.59%  "
# Row 16
$ws.Range("D16").Value = "1.778.79"
$ws.Range("E16").Value = "  +1.46%  "
# Row 17
$ws.Range("D17").Value = "'93.14"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.11%  "
# Row 18
$ws.Range("D18").Value = "'0.00001061"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.19%  "
# Row 19
$ws.Range("D19").Value = "'0.06449"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.64%  "
# Row 20
$ws.Range("D20").Value = "'1.004"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.37%  "
# Row 21
$ws.Range("D21").Value = "'17.22"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.55%  "
# Row 22
$ws.Range("D22").Value = "'5.812"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.38%  "
# Row 23
$ws.Range("D23").Value = "27.908.82"
$ws.Range("E23").Value = "  +0.99%  "
# Row 24
$ws.Range("D24").Value = "'11.33"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.21%  "
# Row 25
$ws.Range("D25").Value = "'2.110"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.28%  "
# Row 26
$ws.Range("D26").Value = "'163.24"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.38%  "
# Row 27
$ws.Range("D27").Value = "'20.30"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.55%  "
# Row 28
$ws.Range("D28").Value = "1.984.49"
$ws.Range("E28").Value = "  +1.53%  "
# Row 29
$ws.Range("D29").Value = "'2.208"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.12%  "
# Row 30
$ws.Range("D30").Value = "'125.67"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.33%  "
# Row 31
$ws.Range("D31").Value = "'1.102"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.44%  "
# Row 32
$ws.Range("D32").Value = "'0.09119"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.04%  "
# Row 33
$ws.Range("D33").Value = "'5.566"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.84%  "
# Row 34
$ws.Range("D34").Value = "'3.644"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.53%  "
# Row 35
$ws.Range("D35").Value = "'11.89"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.11%  "
# Row 36
$ws.Range("D36").Value = "'0.02293"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.19%  "
# Row 37
$ws.Range("D37").Value = "'0.06110"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.31%  "
# Row 38
$ws.Range("D38").Value = "'0.2092"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.59%  "
# Row 39
$ws.Range("D39").Value = "'0.6342"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.01%  "
# Row 40
$ws.Range("D40").Value = "'4.960"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.12%  "
# Row 41
$ws.Range("D41").Value = "'1.185"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.25%  "
# Row 42
$ws.Range("D42").Value = "'1.401"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.27%  "
# Row 43
$ws.Range("D43").Value = "'7.928"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.29%  "
# Row 44
$ws.Range("D44").Value = "'13.33"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.58%  "
# Row 45
$ws.Range("B45").Value = "PancakeSwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D45").Value = "'3.733"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.54%  "
# Row 46
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "'0.5908"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.31%  "
# Row 47
$ws.Range("D47").Value = "'122.77"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.14%  "
# Row 48
$ws.Range("D48").Value = "'1.958"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.37%  "
# Row 49
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.06921"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.96%  "
# Row 50
$ws.Range("B50").Value = "EOS"
$ws.Range("C50").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D50").Value = "'1.138"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.89%  "
# Row 51
$ws.Range("D51").Value = "'72.81"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.02%  "
